# Applies the OOXML diff:
#  - MTHousing: rows 3,6,9,12,15,18,21,24 cols I:K -> "-1" (was "-1.0"/"nan")
#  - MTHousing: append two new data rows (27, 28)
#  - EZLineHousing, REDAHousing, OtherHousing: row 2 cols I:K -> "-1"
#    (was 4860 / 4866 / 4854)
#
# Helper: Excel.Range.Value auto-detects numeric-looking strings ("-1",
# "410.0", ...) and stores them as Number cells. The source file stores
# these as plain Text, so we briefly force a Text number format, assign
# the value, then restore the cell style to "Normal" (General format /
# default style) so no stray formatting is left behind.
function Set-TextValue {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MTHousing
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("MTHousing")

# rows with stray "-1.0"/"nan" markers -> normalize to "-1"
$fixRows = @(3, 6, 9, 12, 15, 18, 21, 24)
foreach ($r in $fixRows) {
    foreach ($col in 9, 10, 11) {
        Set-TextValue $ws $r $col "-1"
    }
}

# Append two new rows, copying row 26's formatting (keeps column A's
# numeric style and columns B:K as text) and then overwriting values.
$ws.Range("A26:K26").Copy($ws.Range("A27:K27"))
$ws.Range("A26:K26").Copy($ws.Range("A28:K28"))

# Row 27
$ws.Cells.Item(27, 1).Value = 25
Set-TextValue $ws 27 2 "TPS-Line"
Set-TextValue $ws 27 3 "2A"
Set-TextValue $ws 27 4 "CR"
Set-TextValue $ws 27 5 "3"
Set-TextValue $ws 27 6 "3000"
Set-TextValue $ws 27 7 "2"
Set-TextValue $ws 27 8 "3"
Set-TextValue $ws 27 9 "-1"
Set-TextValue $ws 27 10 "-1"
Set-TextValue $ws 27 11 "-1"

# Row 28
$ws.Cells.Item(28, 1).Value = 26
Set-TextValue $ws 28 2 "TPS-Line"
Set-TextValue $ws 28 3 "2A"
Set-TextValue $ws 28 4 "FL"
Set-TextValue $ws 28 5 "0.5"
Set-TextValue $ws 28 6 "500"
Set-TextValue $ws 28 7 "1"
Set-TextValue $ws 28 8 "1"
Set-TextValue $ws 28 9 "410.0"
Set-TextValue $ws 28 10 "413.2"
Set-TextValue $ws 28 11 "406.8"

# ---------------------------------------------------------------------
# EZLineHousing / REDAHousing / OtherHousing: row 2, cols I:K -> "-1"
# ---------------------------------------------------------------------
$otherSheets = @("EZLineHousing", "REDAHousing", "OtherHousing")
foreach ($sheetName in $otherSheets) {
    $ws2 = $wb.Worksheets.Item($sheetName)
    Set-TextValue $ws2 2 9 "-1"
    Set-TextValue $ws2 2 10 "-1"
    Set-TextValue $ws2 2 11 "-1"
}
